# 0.2.3 added excel products
# Adds a new "product_discount" column (K) with default (0) values for the
# existing products, fills in the previously-empty "id" column (A) for the
# existing rows, and highlights the new discount cells in yellow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty "id" column (A) to mirror "product_code" (B)
$ws.Range("A2").Value = 2001
$ws.Range("A3").Value = 2002

# New column K: header + data for the new "product_discount" field
$ws.Range("K1").Value = "product_discount"
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0

# Match the font size used by the rest of the data rows
$ws.Range("K2:K3").Font.Size = 12

# Highlight the new discount values in yellow
$ws.Range("K2:K3").Interior.Color = 65535

# Size the new column to fit its header text
$ws.Columns.Item(11).ColumnWidth = 13.5

# Reflect the rows that were edited as the active selection
[void]$ws.Range("A2:A3").EntireRow.Select()
